# Update weekly excess mortality data and plots
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("excess_mortality_provinces")

# --- Revised weekly figures for already-present weeks (rows 117-147) ---
$ws.Range("X117").Value = 674
$ws.Range("X118").Value = 722
$ws.Range("W119").Value = 567
$ws.Range("AA121").Value = 305
$ws.Range("X131").Value = 707
$ws.Range("Z131").Value = 453
$ws.Range("V135").Value = 195
$ws.Range("W135").Value = 510
$ws.Range("X136").Value = 601
$ws.Range("W138").Value = 447
$ws.Range("W139").Value = 502
$ws.Range("X139").Value = 606
$ws.Range("W140").Value = 444
$ws.Range("X140").Value = 595
$ws.Range("Z140").Value = 489
$ws.Range("U141").Value = 395
$ws.Range("X141").Value = 599
$ws.Range("AA141").Value = 239
$ws.Range("W142").Value = 457
$ws.Range("X143").Value = 604
$ws.Range("AA143").Value = 224
$ws.Range("R144").Value = 109
$ws.Range("U144").Value = 377
$ws.Range("W144").Value = 485
$ws.Range("X144").Value = 606
$ws.Range("AA144").Value = 221
$ws.Range("S145").Value = 224
$ws.Range("W145").Value = 480
$ws.Range("X145").Value = 654
$ws.Range("P146").Value = 129
$ws.Range("S146").Value = 244
$ws.Range("U146").Value = 447
$ws.Range("W146").Value = 495
$ws.Range("X146").Value = 658
$ws.Range("Z146").Value = 500
$ws.Range("AA146").Value = 222
$ws.Range("P147").Value = 106
$ws.Range("Q147").Value = 124
$ws.Range("R147").Value = 106
$ws.Range("S147").Value = 231
$ws.Range("T147").Value = 49
$ws.Range("U147").Value = 458
$ws.Range("V147").Value = 212
$ws.Range("W147").Value = 468
$ws.Range("X147").Value = 680
$ws.Range("Z147").Value = 480
$ws.Range("AA147").Value = 255

# --- New row 148: 2022 week 42 ---
$ws.Range("N148").Value = 2022
$ws.Range("O148").Value = 42
$ws.Range("P148").Value = 116
$ws.Range("Q148").Value = 161
$ws.Range("R148").Value = 109
$ws.Range("S148").Value = 255
$ws.Range("T148").Value = 65
$ws.Range("U148").Value = 429
$ws.Range("V148").Value = 249
$ws.Range("W148").Value = 455
$ws.Range("X148").Value = 634
$ws.Range("Y148").Value = 79
$ws.Range("Z148").Value = 527
$ws.Range("AA148").Value = 240

$ws.Range("AC148").Value = 2022
$ws.Range("AD148").Value = 42

$ws.Range("AE148").Formula = "=ROUND((P148-B148)/B148*100,2)"
$ws.Range("AF148").Formula = "=ROUND((Q148-C148)/C148*100,2)"
$ws.Range("AG148").Formula = "=ROUND((R148-D148)/D148*100,2)"
$ws.Range("AH148").Formula = "=ROUND((S148-E148)/E148*100,2)"
$ws.Range("AI148").Formula = "=ROUND((T148-F148)/F148*100,2)"
$ws.Range("AJ148").Formula = "=ROUND((U148-G148)/G148*100,2)"
$ws.Range("AK148").Formula = "=ROUND((V148-H148)/H148*100,2)"
$ws.Range("AL148").Formula = "=ROUND((W148-I148)/I148*100,2)"
$ws.Range("AM148").Formula = "=ROUND((X148-J148)/J148*100,2)"
$ws.Range("AN148").Formula = "=ROUND((Y148-K148)/K148*100,2)"
$ws.Range("AO148").Formula = "=ROUND((Z148-L148)/L148*100,2)"
$ws.Range("AP148").Formula = "=ROUND((AA148-M148)/M148*100,2)"

# --- View state (sheet selection / scroll position) ---
$ws.Application.ActiveWindow.ScrollRow = 106
$ws.Range("AG141").Select() | Out-Null
